$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counter cell C3: keep as literal text "21.0" (not auto-converted to number)
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "21.0"
$ws.Range("C3").ClearFormats()

# Widen column B to fit content
$ws.Range("B1").EntireColumn.ColumnWidth = 11.67

# Update the selection to just C2
[void]$ws.Range("C2").Select()
